# Scheduled-runner price refresh: rewrites the market-price-derived
# columns (H:N -- currentAveragePrice[NQ/HQ], LevePrice[NQ/HQ],
# LeveProfit[NQ/HQ]) for the leves whose Universalis quotes changed,
# one worksheet per crafting job (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
# Rows whose market data collapsed to all-zero (no NQ/HQ listings)
# drop their LeveProfit cells entirely, matching how the rest of the
# sheet represents "no data" rows.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 299.75
$ws.Range("I39").Value = 66.333336
$ws.Range("J39").Value = 1000
$ws.Range("K39").Value = 199.000008
$ws.Range("L39").Value = 3000
$ws.Range("M39").Value = 96.99999199999999
$ws.Range("N39").Value = -3592

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 45193.723
$ws.Range("I106").Value = 1745.5555
$ws.Range("J106").Value = 88641.89
$ws.Range("K106").Value = 1745.5555
$ws.Range("L106").Value = 88641.89
$ws.Range("M106").Value = -1114.5555
$ws.Range("N106").Value = -89903.89

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1891.7916
$ws.Range("J112").Value = 1957.2858
$ws.Range("L112").Value = 5871.857400000001
$ws.Range("N112").Value = -8087.857400000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 16900.55
$ws.Range("I132").Value = 2452.8845
$ws.Range("J132").Value = 110810.375
$ws.Range("K132").Value = 7358.6535
$ws.Range("L132").Value = 332431.125
$ws.Range("M132").Value = -4828.6535
$ws.Range("N132").Value = -337491.125

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 4105.6304
$ws.Range("I137").Value = 1352.5
$ws.Range("J137").Value = 8388.277
$ws.Range("K137").Value = 4057.5
$ws.Range("L137").Value = 25164.831
$ws.Range("M137").Value = -1507.5
$ws.Range("N137").Value = -30264.831

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1627.25
$ws.Range("I141").Value = 767.87805
$ws.Range("J141").Value = 6660.7144
$ws.Range("K141").Value = 2303.63415
$ws.Range("L141").Value = 19982.1432
$ws.Range("M141").Value = 2876.36585
$ws.Range("N141").Value = -30342.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11087.804
$ws.Range("I32").Value = 10026.075
$ws.Range("K32").Value = 10026.075
$ws.Range("M32").Value = -9739.075000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1394.6041
$ws.Range("I61").Value = 1193.909
$ws.Range("J61").Value = 3602.25
$ws.Range("K61").Value = 1193.909
$ws.Range("L61").Value = 3602.25
$ws.Range("M61").Value = -981.9090000000001
$ws.Range("N61").Value = -4026.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1512.6111
$ws.Range("I74").Value = 1551.4865
$ws.Range("J74").Value = 1428
$ws.Range("K74").Value = 1551.4865
$ws.Range("L74").Value = 1428
$ws.Range("M74").Value = -677.4865
$ws.Range("N74").Value = -3176

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1512.6111
$ws.Range("I77").Value = 1551.4865
$ws.Range("J77").Value = 1428
$ws.Range("K77").Value = 7757.4325
$ws.Range("L77").Value = 7140
$ws.Range("M77").Value = -3389.4325
$ws.Range("N77").Value = -15876

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1394.6041
$ws.Range("I136").Value = 1193.909
$ws.Range("J136").Value = 3602.25
$ws.Range("K136").Value = 3581.727
$ws.Range("L136").Value = 10806.75
$ws.Range("M136").Value = -1031.727
$ws.Range("N136").Value = -15906.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3350.2986
$ws.Range("I134").Value = 1900.7646
$ws.Range("J134").Value = 3843.14
$ws.Range("K134").Value = 5702.293799999999
$ws.Range("L134").Value = 11529.42
$ws.Range("M134").Value = -3167.293799999999
$ws.Range("N134").Value = -16599.42

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2312.818
$ws.Range("I22").Value = 410.25
$ws.Range("J22").Value = 3400
$ws.Range("K22").Value = 410.25
$ws.Range("L22").Value = 3400
$ws.Range("M22").Value = -60.25
$ws.Range("N22").Value = -4100

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2866
$ws.Range("I31").Value = 1279.7037
$ws.Range("J31").Value = 3452.7124
$ws.Range("K31").Value = 1279.7037
$ws.Range("L31").Value = 3452.7124
$ws.Range("M31").Value = -984.7037
$ws.Range("N31").Value = -4042.7124

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2866
$ws.Range("I34").Value = 1279.7037
$ws.Range("J34").Value = 3452.7124
$ws.Range("K34").Value = 1279.7037
$ws.Range("L34").Value = 3452.7124
$ws.Range("M34").Value = -1077.7037
$ws.Range("N34").Value = -3856.7124

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1308.3077
$ws.Range("I58").Value = 993.55817
$ws.Range("J58").Value = 2812.111
$ws.Range("K58").Value = 993.55817
$ws.Range("L58").Value = 2812.111
$ws.Range("M58").Value = -790.55817
$ws.Range("N58").Value = -3218.111

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 28486.385
$ws.Range("I132").Value = 1347.341
$ws.Range("J132").Value = 177751.12
$ws.Range("K132").Value = 4042.023
$ws.Range("L132").Value = 533253.36
$ws.Range("M132").Value = -1512.023
$ws.Range("N132").Value = -538313.36

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2128.4644
$ws.Range("I134").Value = 1643.7646
$ws.Range("J134").Value = 2877.5454
$ws.Range("K134").Value = 4931.293799999999
$ws.Range("L134").Value = 8632.636200000001
$ws.Range("M134").Value = -2396.293799999999
$ws.Range("N134").Value = -13702.6362

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1308.3077
$ws.Range("I136").Value = 993.55817
$ws.Range("J136").Value = 2812.111
$ws.Range("K136").Value = 2980.67451
$ws.Range("L136").Value = 8436.332999999999
$ws.Range("M136").Value = -430.6745099999998
$ws.Range("N136").Value = -13536.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1278.5714
$ws.Range("I4").Value = 87.25
$ws.Range("K4").Value = 261.75
$ws.Range("M4").Value = -149.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1014.61536
$ws.Range("I117").Value = 965
$ws.Range("J117").Value = 1057.1428
$ws.Range("K117").Value = 2895
$ws.Range("L117").Value = 3171.4284
$ws.Range("M117").Value = 547
$ws.Range("N117").Value = -10055.4284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2032.2742
$ws.Range("I132").Value = 1368.2632
$ws.Range("J132").Value = 3083.625
$ws.Range("K132").Value = 4104.7896
$ws.Range("L132").Value = 9250.875
$ws.Range("M132").Value = -1574.7896
$ws.Range("N132").Value = -14310.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 78839.38
$ws.Range("I122").Value = 101690.7
$ws.Range("J122").Value = 2668.3333
$ws.Range("K122").Value = 305072.1
$ws.Range("L122").Value = 8004.999899999999
$ws.Range("M122").Value = -302622.1
$ws.Range("N122").Value = -12904.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2799.2708
$ws.Range("I132").Value = 1984.0938
$ws.Range("J132").Value = 4429.625
$ws.Range("K132").Value = 5952.2814
$ws.Range("L132").Value = 13288.875
$ws.Range("M132").Value = -3422.2814
$ws.Range("N132").Value = -18348.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1449.6279
$ws.Range("I136").Value = 1144.5625
$ws.Range("J136").Value = 2337.0908
$ws.Range("K136").Value = 3433.6875
$ws.Range("L136").Value = 7011.2724
$ws.Range("M136").Value = -883.6875
$ws.Range("N136").Value = -12111.2724

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 14308.5
$ws.Range("J74").Value = 13250.2
$ws.Range("L74").Value = 13250.2
$ws.Range("N74").Value = -15122.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H77").Value = 14308.5
$ws.Range("J77").Value = 13250.2
$ws.Range("L77").Value = 39750.60000000001
$ws.Range("N77").Value = -49110.60000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1842.1428
$ws.Range("I132").Value = 1388.2759
$ws.Range("J132").Value = 2854.6155
$ws.Range("K132").Value = 4164.8277
$ws.Range("L132").Value = 8563.8465
$ws.Range("M132").Value = -1634.8277
$ws.Range("N132").Value = -13623.8465

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 25477.586
$ws.Range("I136").Value = 39158.848
$ws.Range("J136").Value = 1763.4
$ws.Range("K136").Value = 117476.544
$ws.Range("L136").Value = 5290.200000000001
$ws.Range("M136").Value = -114926.544
$ws.Range("N136").Value = -10390.2
